$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns P and Q ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Copy the existing header formatting (bold font, thin border, centered
# horizontal / top vertical alignment) from O1 onto the two new header
# cells so they match the rest of the header row.
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Data rows 2-25: update I/K/M/O and add P/Q ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
